$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3547.5557
$ws.Range("I113").Value = 2660.4546
$ws.Range("J113").Value = 4157.4375
$ws.Range("K113").Value = 2660.4546
$ws.Range("L113").Value = 4157.4375
$ws.Range("M113").Value = 593.5454
$ws.Range("N113").Value = -10665.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4026178
$ws.Range("I32").Value = 5783.722
$ws.Range("J32").Value = 30341486
$ws.Range("K32").Value = 5783.722
$ws.Range("L32").Value = 30341486
$ws.Range("M32").Value = -5496.722
$ws.Range("N32").Value = -30342060

$ws.Range("H110").Value = 3527.25
$ws.Range("I110").Value = 2494.2222
$ws.Range("J110").Value = 4855.4287
$ws.Range("K110").Value = 2494.2222
$ws.Range("L110").Value = 4855.4287
$ws.Range("M110").Value = -449.2222000000002
$ws.Range("N110").Value = -8945.4287

$ws.Range("H132").Value = 2554829.5
$ws.Range("I132").Value = 1454.5161
$ws.Range("K132").Value = 4363.5483
$ws.Range("M132").Value = -1833.5483

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1598.2
$ws.Range("I20").Value = 1282.8
$ws.Range("J20").Value = 2386.7
$ws.Range("K20").Value = 1282.8
$ws.Range("L20").Value = 2386.7
$ws.Range("M20").Value = -1035.8
$ws.Range("N20").Value = -2880.7

$ws.Range("H86").Value = 1911.125
$ws.Range("I86").Value = 1884.1428
$ws.Range("K86").Value = 1884.1428
$ws.Range("M86").Value = -761.1428000000001

$ws.Range("H89").Value = 1911.125
$ws.Range("I89").Value = 1884.1428
$ws.Range("K89").Value = 9420.714
$ws.Range("M89").Value = -3804.714

$ws.Range("H99").Value = 2015.6428
$ws.Range("I99").Value = 1834.4706
$ws.Range("J99").Value = 2295.6365
$ws.Range("K99").Value = 1834.4706
$ws.Range("L99").Value = 2295.6365
$ws.Range("M99").Value = -336.4706000000001
$ws.Range("N99").Value = -5291.636500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4568009.5
$ws.Range("I31").Value = 1492.2903
$ws.Range("J31").Value = 7938534
$ws.Range("K31").Value = 1492.2903
$ws.Range("L31").Value = 7938534
$ws.Range("M31").Value = -1197.2903
$ws.Range("N31").Value = -7939124

$ws.Range("H34").Value = 4568009.5
$ws.Range("I34").Value = 1492.2903
$ws.Range("J34").Value = 7938534
$ws.Range("K34").Value = 1492.2903
$ws.Range("L34").Value = 7938534
$ws.Range("M34").Value = -1290.2903
$ws.Range("N34").Value = -7938938

$ws.Range("H99").Value = 7695502.5
$ws.Range("I99").Value = 1941.75
$ws.Range("J99").Value = 20005200
$ws.Range("K99").Value = 1941.75
$ws.Range("L99").Value = 20005200
$ws.Range("M99").Value = -443.75
$ws.Range("N99").Value = -20008196

$ws.Range("H126").Value = 7695502.5
$ws.Range("I126").Value = 1941.75
$ws.Range("J126").Value = 20005200
$ws.Range("K126").Value = 5825.25
$ws.Range("L126").Value = 60015600
$ws.Range("M126").Value = -3355.25
$ws.Range("N126").Value = -60020540

$ws.Range("H132").Value = 2824.4666
$ws.Range("I132").Value = 2184.2
$ws.Range("K132").Value = 6552.599999999999
$ws.Range("M132").Value = -4022.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 700512.5600000001
$ws.Range("I68").Value = 942.5714
$ws.Range("J68").Value = 1034398.2
$ws.Range("K68").Value = 2827.7142
$ws.Range("L68").Value = 3103194.6
$ws.Range("M68").Value = -2016.7142
$ws.Range("N68").Value = -3104816.6

$ws.Range("H71").Value = 700512.5600000001
$ws.Range("I71").Value = 942.5714
$ws.Range("J71").Value = 1034398.2
$ws.Range("K71").Value = 8483.142600000001
$ws.Range("L71").Value = 9309583.799999999
$ws.Range("M71").Value = -4427.142600000001
$ws.Range("N71").Value = -9317695.799999999

$ws.Range("H80").Value = 2483.8462
$ws.Range("J80").Value = 2574.1667
$ws.Range("L80").Value = 7722.500100000001
$ws.Range("N80").Value = -9594.500100000001

$ws.Range("H83").Value = 2483.8462
$ws.Range("J83").Value = 2574.1667
$ws.Range("L83").Value = 23167.5003
$ws.Range("N83").Value = -32527.5003

$ws.Range("H113").Value = 505.69232
$ws.Range("I113").Value = 487.66666
$ws.Range("J113").Value = 521.1429000000001
$ws.Range("K113").Value = 1462.99998
$ws.Range("L113").Value = 1563.4287
$ws.Range("M113").Value = 707.0000199999999
$ws.Range("N113").Value = -5903.4287

$ws.Range("H122").Value = 3395.3093
$ws.Range("I122").Value = 295
$ws.Range("J122").Value = 3712.3865
$ws.Range("K122").Value = 2655
$ws.Range("L122").Value = 33411.4785
$ws.Range("M122").Value = -205
$ws.Range("N122").Value = -38311.4785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2553752
$ws.Range("I102").Value = 4466201
$ws.Range("J102").Value = 3819.9167
$ws.Range("K102").Value = 4466201
$ws.Range("L102").Value = 3819.9167
$ws.Range("M102").Value = -4464579
$ws.Range("N102").Value = -7063.9167

$ws.Range("H126").Value = 20835682
$ws.Range("I126").Value = 50001024
$ws.Range("J126").Value = 3296.2856
$ws.Range("K126").Value = 150003072
$ws.Range("L126").Value = 9888.856800000001
$ws.Range("M126").Value = -150000602
$ws.Range("N126").Value = -14828.8568

$ws.Range("H132").Value = 2285.0356
$ws.Range("I132").Value = 2634.5
$ws.Range("K132").Value = 7903.5
$ws.Range("M132").Value = -5373.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3232.9412
$ws.Range("I68").Value = 2873
$ws.Range("J68").Value = 3747.1428
$ws.Range("K68").Value = 2873
$ws.Range("L68").Value = 3747.1428
$ws.Range("M68").Value = -2124
$ws.Range("N68").Value = -5245.1428

$ws.Range("H71").Value = 3232.9412
$ws.Range("I71").Value = 2873
$ws.Range("J71").Value = 3747.1428
$ws.Range("K71").Value = 14365
$ws.Range("L71").Value = 18735.714
$ws.Range("M71").Value = -10621
$ws.Range("N71").Value = -26223.714

$ws.Range("H82").Value = 3259.5
$ws.Range("I82").Value = 2650
$ws.Range("J82").Value = 3564.25
$ws.Range("K82").Value = 2650
$ws.Range("L82").Value = 3564.25
$ws.Range("M82").Value = -2289
$ws.Range("N82").Value = -4286.25

$ws.Range("H85").Value = 3259.5
$ws.Range("I85").Value = 2650
$ws.Range("J85").Value = 3564.25
$ws.Range("K85").Value = 2650
$ws.Range("L85").Value = 3564.25
$ws.Range("M85").Value = -1402
$ws.Range("N85").Value = -6060.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7401.8184
$ws.Range("I96").Value = 2990
$ws.Range("J96").Value = 8382.223
$ws.Range("K96").Value = 2990
$ws.Range("L96").Value = 8382.223
$ws.Range("M96").Value = -1617
$ws.Range("N96").Value = -11128.223

$ws.Range("H122").Value = 9290.35
$ws.Range("I122").Value = 11762.462
$ws.Range("J122").Value = 4699.2856
$ws.Range("K122").Value = 35287.386
$ws.Range("L122").Value = 14097.8568
$ws.Range("M122").Value = -32837.386
$ws.Range("N122").Value = -18997.8568
